$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126-200 down to 127-201.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new data record.
$ws.Cells.Item(126, 1).Value = 10
$ws.Cells.Item(126, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(126, 3).Value = "La Araucanía"
$ws.Cells.Item(126, 4).Value = 44452
$ws.Cells.Item(126, 5).Value = 9
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100108
$ws.Cells.Item(126, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(126, 9).Value = 100108002
$ws.Cells.Item(126, 10).Value = "Mango"
$ws.Cells.Item(126, 11).Value = "Sin especificar"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 385
$ws.Cells.Item(126, 14).Value = 9000
$ws.Cells.Item(126, 15).Value = 9000
$ws.Cells.Item(126, 16).Value = 9000
$ws.Cells.Item(126, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(126, 18).Value = "Brasil"
$ws.Cells.Item(126, 19).Value = 2250
$ws.Cells.Item(126, 20).Value = 4
